$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) contains values such as "1.009" or "0.2200" that look
# numeric but must stay literal text (matching the original inline-string cells,
# including formatting such as trailing zeros). Force the whole edited range to
# Text format first so Excel does not auto-convert these assignments to numbers,
# then restore the default "Normal" style afterwards so no stray number format is
# left behind on the cells.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.680.47"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.883.62"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  -1.87%  "
$ws.Range("D5").Value = "316.35"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  -1.77%  "
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("D8").Value = "0.3957"
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("D9").Value = "0.08445"
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("D10").Value = "1.117"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "41.81"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "6.308"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.886.02"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").Value = "20.57"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "7.304"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.00001110"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "91.61"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "0.06741"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "17.80"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.008"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "5.991"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("B23").Value = "WrappedBTC"
$ws.Range("C23").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D23").Value = "28.713.37"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "11.17"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.256"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.100.39"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "161.18"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "20.78"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "2.391"
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "127.27"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.1055"
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "1.054"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.826"
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.629"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "0.02475"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.06557"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "0.2202"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "8.971"
$ws.Range("E38").Value = "  -4.44%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "1.265"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.203"
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").Value = "5.114"
$ws.Range("E41").Value = "  +2.63%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.6484"
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "11.23"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "1.008"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.6089"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "13.04"
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.708"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "2.046"
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "1.225"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "122.78"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("B51").Value = "WEMIXTOKEN"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "1.198"
$ws.Range("E51").Value = "  -7.82%  "

# Restore default styling on the whole edited range.
$editRange.Style = "Normal"
